# Add data for 2021-12-06: extend "November (through 11-27)" data to
# "November (through 11-28)" — updates the workbook/sheet title, the row
# label, and the 2016-2020 + 2021 columns for the November row (13) and
# the Total row (14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and, implicitly, the workbook's sheet tab name).
$ws.Name = "Through 2021-11-28"

# Update the row label for November.
$ws.Range("A13").Value = "November (through 11-28)"

# Row 13 ("November (through 11-28)") updates.
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 65
$ws.Range("G13").Value = 0.1096

$ws.Range("I13").Value = 99
$ws.Range("J13").Value = 0.0198

$ws.Range("L13").Value = 55
$ws.Range("M13").Value = 0.127

$ws.Range("O13").Value = 42
$ws.Range("P13").Value = 0.125

$ws.Range("R13").Value = 189
$ws.Range("S13").Value = 0.0503

$ws.Range("U13").Value = 186
$ws.Range("V13").Value = 0.0211

# Row 14 ("Total") updates.
$ws.Range("E14").Value = 60
$ws.Range("F14").Value = 499
$ws.Range("G14").Value = 0.1073

$ws.Range("I14").Value = 748
$ws.Range("J14").Value = 0.0777

$ws.Range("L14").Value = 604
$ws.Range("M14").Value = 0.1091

$ws.Range("O14").Value = 476
$ws.Range("P14").Value = 0.1019

$ws.Range("R14").Value = 1193
$ws.Range("S14").Value = 0.0502

$ws.Range("U14").Value = 1537
$ws.Range("V14").Value = 0.0588
